$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (Q) mirroring the existing "2019" column (P):
# header cell Q4 (year) and data cell Q5 (value), copying P's formatting.

$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("Q5").Value = 90.6

$excel.CutCopyMode = $false

# Move/record the active selection as it was left in the authored workbook.
$ws.Range("P12").Select() | Out-Null
